$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.305.06'
$ws.Range("E2").Value = '  +2.95%  '
$ws.Range("D3").Value = '2.310.36'
$ws.Range("E3").Value = '  +2.72%  '
$ws.Range("D4").Value = "'" + '0.999'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'" + '309.05'
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("D6").Value = "'" + '104.95'
$ws.Range("E6").Value = '  +8.50%  '
$ws.Range("E7").Value = '  +0.51%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = "'" + '0.519'
$ws.Range("E9").Value = '  +6.55%  '
$ws.Range("D10").Value = "'" + '35.92'
$ws.Range("E10").Value = '  +3.10%  '
$ws.Range("D11").Value = "'" + '52.68'
$ws.Range("E11").Value = '  +2.16%  '
$ws.Range("D12").Value = "'" + '0.0809'
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("E13").Value = '  -1.30%  '
$ws.Range("D14").Value = "'" + '6.96'
$ws.Range("E14").Value = '  +2.80%  '
$ws.Range("D15").Value = '2.668.55'
$ws.Range("E15").Value = '  +2.71%  '
$ws.Range("D16").Value = "'" + '15.17'
$ws.Range("E16").Value = '  +5.24%  '
$ws.Range("D17").Value = '2.302.99'
$ws.Range("E17").Value = '  +2.57%  '
$ws.Range("D18").Value = "'" + '0.801'
$ws.Range("E18").Value = '  +2.91%  '
$ws.Range("D19").Value = '43.232.34'
$ws.Range("E19").Value = '  +3.01%  '
$ws.Range("D20").Value = "'" + '11.98'
$ws.Range("E20").Value = '  -1.27%  '
$ws.Range("D21").Value = '0.0₃0922'
$ws.Range("E21").Value = '  +2.28%  '
$ws.Range("D22").Value = "'" + '6.20'
$ws.Range("E22").Value = '  +4.78%  '
$ws.Range("D23").Value = "'" + '67.89'
$ws.Range("E23").Value = '  +1.05%  '
$ws.Range("D24").Value = "'" + '240.21'
$ws.Range("E24").Value = '  +2.14%  '
$ws.Range("D25").Value = "'" + '2.02'
$ws.Range("E25").Value = '  +3.92%  '
$ws.Range("E26").Value = '  +1.62%  '
$ws.Range("E27").Value = '  -1.45%  '
$ws.Range("D28").Value = "'" + '24.75'
$ws.Range("E28").Value = '  +6.43%  '
$ws.Range("D29").Value = "'" + '2.22'
$ws.Range("E29").Value = '  +4.77%  '
$ws.Range("D30").Value = "'" + '36.18'
$ws.Range("E30").Value = '  -4.50%  '
$ws.Range("D31").Value = "'" + '9.58'
$ws.Range("E31").Value = '  +1.22%  '
$ws.Range("D32").Value = "'" + '162.33'
$ws.Range("E32").Value = '  -3.01%  '
$ws.Range("D33").Value = "'" + '5.24'
$ws.Range("E33").Value = '  +1.81%  '
$ws.Range("D34").Value = "'" + '0.998'
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("D35").Value = "'" + '18.36'
$ws.Range("E35").Value = '  +5.20%  '
$ws.Range("E36").Value = '  +6.36%  '
$ws.Range("D37").Value = "'" + '0.0734'
$ws.Range("E37").Value = '  +2.09%  '
$ws.Range("D38").Value = "'" + '3.00'
$ws.Range("E38").Value = '  -2.41%  '
$ws.Range("D39").Value = "'" + '4.53'
$ws.Range("E39").Value = '  +11.82%  '
$ws.Range("D40").Value = "'" + '0.106'
$ws.Range("E40").Value = '  +3.26%  '
$ws.Range("D41").Value = "'" + '1.86'
$ws.Range("E41").Value = '  +4.10%  '
$ws.Range("E42").Value = '  +0.51%  '
$ws.Range("D43").Value = "'" + '2.46'
$ws.Range("E43").Value = '  +13.05%  '
$ws.Range("D44").Value = "'" + '0.0289'
$ws.Range("E44").Value = '  +2.91%  '
$ws.Range("D45").Value = '1.966.91'
$ws.Range("E45").Value = '  +1.46%  '
$ws.Range("D46").Value = "'" + '18.78'
$ws.Range("E46").Value = '  +1.06%  '
$ws.Range("D47").Value = "'" + '3.07'
$ws.Range("E47").Value = '  +6.14%  '
$ws.Range("D48").Value = "'" + '10.19'
$ws.Range("E48").Value = '  +5.71%  '
$ws.Range("D49").Value = "'" + '57.82'
$ws.Range("E49").Value = '  +7.54%  '
$ws.Range("D50").Value = "'" + '2.94'
$ws.Range("E50").Value = '  +0.74%  '
$ws.Range("E51").Value = '  +8.73%  '

